$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Release version" column (A2:A33) is being bumped from "1.0" to "1.1.0".
# Use a leading apostrophe via .Formula so the cell keeps being stored as
# explicit text (quote-prefixed, right-aligned) instead of the engine
# re-inferring a fresh style for the assigned value.
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 1).Formula = "'1.1.0"
}

# Reset the view: drop the scrolled-down topLeftCell and select column A
# (whole-column selection) instead of the single cell D32.
$ws.Range("A1:A1048576").Select()
